$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "29.688.30"
$ws.Range("E2").Value = "  -2.72%  "

# Row 3
Set-TextValue $ws.Range("D3") "2.094.47"
$ws.Range("E3").Value = "  -1.87%  "

# Row 4
$ws.Range("E4").Value = "  +0.35%  "

# Row 5
Set-TextValue $ws.Range("D5") "343.50"
$ws.Range("E5").Value = "  -2.11%  "

# Row 6
$ws.Range("E6").Value = "  +0.27%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.5161"
$ws.Range("E7").Value = "  -1.75%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.4380"
$ws.Range("E8").Value = "  -3.68%  "

# Row 9
Set-TextValue $ws.Range("D9") "52.76"
$ws.Range("E9").Value = "  -1.48%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.09238"
$ws.Range("E10").Value = "  +1.40%  "

# Row 11
Set-TextValue $ws.Range("D11") "1.163"
$ws.Range("E11").Value = "  -2.46%  "

# Row 12
Set-TextValue $ws.Range("D12") "24.78"
$ws.Range("E12").Value = "  -2.79%  "

# Row 13
Set-TextValue $ws.Range("D13") "2.101.17"
$ws.Range("E13").Value = "  -1.57%  "

# Row 14
Set-TextValue $ws.Range("D14") "8.231"
$ws.Range("E14").Value = "  +0.44%  "

# Row 15
Set-TextValue $ws.Range("D15") "6.756"
$ws.Range("E15").Value = "  -1.85%  "

# Row 16
Set-TextValue $ws.Range("D16") "99.59"
$ws.Range("E16").Value = "  -1.34%  "

# Row 17
Set-TextValue $ws.Range("D17") "0.00001150"
$ws.Range("E17").Value = "  -1.44%  "

# Row 18
$ws.Range("E18").Value = "  +0.30%  "

# Row 19
$ws.Range("E19").Value = "  +1.60%  "

# Row 20
Set-TextValue $ws.Range("D20") "0.06660"
$ws.Range("E20").Value = "  -0.85%  "

# Row 21
$ws.Range("E21").Value = "  +0.27%  "

# Row 22
Set-TextValue $ws.Range("D22") "6.186"
$ws.Range("E22").Value = "  -2.73%  "

# Row 23
Set-TextValue $ws.Range("D23") "29.738.07"
$ws.Range("E23").Value = "  -2.75%  "

# Row 24
$ws.Range("E24").Value = "  -2.63%  "

# Row 26
Set-TextValue $ws.Range("D26") "2.346.15"
$ws.Range("E26").Value = "  -1.67%  "

# Row 27
Set-TextValue $ws.Range("D27") "21.95"
$ws.Range("E27").Value = "  -2.17%  "

# Row 28
Set-TextValue $ws.Range("D28") "2.514"
$ws.Range("E28").Value = "  -3.01%  "

# Row 29
Set-TextValue $ws.Range("D29") "161.11"
$ws.Range("E29").Value = "  -2.15%  "

# Row 30
Set-TextValue $ws.Range("D30") "133.03"
$ws.Range("E30").Value = "  -1.60%  "

# Row 31
Set-TextValue $ws.Range("D31") "1.138"
$ws.Range("E31").Value = "  -6.79%  "

# Row 32
$ws.Range("E32").Value = "  -2.96%  "

# Row 33
Set-TextValue $ws.Range("D33") "1.625"
$ws.Range("E33").Value = "  -5.17%  "

# Row 34
Set-TextValue $ws.Range("D34") "6.180"
$ws.Range("E34").Value = "  -3.32%  "

# Row 35
Set-TextValue $ws.Range("D35") "3.951"
$ws.Range("E35").Value = "  -2.10%  "

# Row 36
Set-TextValue $ws.Range("D36") "6.298"
$ws.Range("E36").Value = "  +2.79%  "

# Row 37
Set-TextValue $ws.Range("D37") "10.18"
$ws.Range("E37").Value = "  -2.59%  "

# Row 38
Set-TextValue $ws.Range("D38") "0.02570"
$ws.Range("E38").Value = "  -2.71%  "

# Row 39
Set-TextValue $ws.Range("D39") "0.7078"
$ws.Range("E39").Value = "  +1.48%  "

# Row 40
Set-TextValue $ws.Range("D40") "0.06709"
$ws.Range("E40").Value = "  -3.58%  "

# Row 41
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D41") "12.43"
$ws.Range("E41").Value = "  -2.27%  "

# Row 42
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D42") "1.329"
$ws.Range("E42").Value = "  +4.71%  "

# Row 43
Set-TextValue $ws.Range("D43") "0.2225"
$ws.Range("E43").Value = "  -5.90%  "

# Row 44
Set-TextValue $ws.Range("D44") "0.6986"
$ws.Range("E44").Value = "  +7.50%  "

# Row 45
Set-TextValue $ws.Range("D45") "14.25"
$ws.Range("E45").Value = "  -3.52%  "

# Row 46
$ws.Range("E46").Value = "  +0.22%  "

# Row 47
Set-TextValue $ws.Range("D47") "2.311"
$ws.Range("E47").Value = "  -1.46%  "

# Row 48
Set-TextValue $ws.Range("D48") "3.622"
$ws.Range("E48").Value = "  -2.15%  "

# Row 49
$ws.Range("E49").Value = "  -5.70%  "

# Row 50
$ws.Range("E50").Value = "  -2.36%  "

# Row 51
Set-TextValue $ws.Range("D51") "81.73"
$ws.Range("E51").Value = "  -2.59%  "
